# Updates the cryptos list with latest price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as text (matching the source workbook's
    # inlineStr cells) instead of letting Excel auto-detect a number/date,
    # while preserving the cell's original style.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "90.952.35"
Set-TextValue $ws.Range("E2") "  +1.42%  "
Set-TextValue $ws.Range("D3") "3.178.97"
Set-TextValue $ws.Range("E3") "  +2.64%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "220.21"
Set-TextValue $ws.Range("E5") "  +2.83%  "
Set-TextValue $ws.Range("D6") "625.73"
Set-TextValue $ws.Range("E6") "  +1.02%  "
Set-TextValue $ws.Range("D7") "1.07"
Set-TextValue $ws.Range("E7") "  +23.77%  "
Set-TextValue $ws.Range("D8") "0.375"
Set-TextValue $ws.Range("E8") "  +1.14%  "
Set-TextValue $ws.Range("E9") "  -0.01%  "
Set-TextValue $ws.Range("D10") "3.178.00"
Set-TextValue $ws.Range("E10") "  +2.71%  "
Set-TextValue $ws.Range("D11") "0.757"
Set-TextValue $ws.Range("E11") "  +14.43%  "
Set-TextValue $ws.Range("D12") "0.198"
Set-TextValue $ws.Range("E12") "  +5.10%  "
Set-TextValue $ws.Range("E13") "  +3.94%  "
Set-TextValue $ws.Range("D14") "35.37"
Set-TextValue $ws.Range("E14") "  +8.74%  "
Set-TextValue $ws.Range("D15") "5.58"
Set-TextValue $ws.Range("E15") "  +3.89%  "
Set-TextValue $ws.Range("D16") "90.438.68"
Set-TextValue $ws.Range("E16") "  +0.65%  "
Set-TextValue $ws.Range("D17") "3.769.81"
Set-TextValue $ws.Range("E17") "  +2.44%  "
Set-TextValue $ws.Range("D18") "3.174.39"
Set-TextValue $ws.Range("E18") "  +2.64%  "
Set-TextValue $ws.Range("D19") "3.78"
Set-TextValue $ws.Range("E19") "  +11.30%  "
Set-TextValue $ws.Range("D20") "0.0000220"
Set-TextValue $ws.Range("E20") "  +2.49%  "
Set-TextValue $ws.Range("D21") "14.39"
Set-TextValue $ws.Range("E21") "  +6.21%  "
Set-TextValue $ws.Range("D22") "443.06"
Set-TextValue $ws.Range("E22") "  +2.49%  "
Set-TextValue $ws.Range("D23") "9.00"
Set-TextValue $ws.Range("E23") "  +8.78%  "
Set-TextValue $ws.Range("D24") "5.22"
Set-TextValue $ws.Range("E24") "  +4.03%  "
Set-TextValue $ws.Range("D25") "6.07"
Set-TextValue $ws.Range("E25") "  +9.08%  "
Set-TextValue $ws.Range("D26") "12.40"
Set-TextValue $ws.Range("E26") "  +1.92%  "
Set-TextValue $ws.Range("D27") "87.25"
Set-TextValue $ws.Range("E27") "  +1.20%  "
Set-TextValue $ws.Range("D28") "3.351.48"
Set-TextValue $ws.Range("E28") "  +2.01%  "
Set-TextValue $ws.Range("E29") "  -0.06%  "
Set-TextValue $ws.Range("D30") "0.164"
Set-TextValue $ws.Range("E30") "  +0.60%  "
Set-TextValue $ws.Range("D31") "9.30"
Set-TextValue $ws.Range("E31") "  +14.19%  "
Set-TextValue $ws.Range("E32") "  -8.95%  "
Set-TextValue $ws.Range("D33") "528.80"
Set-TextValue $ws.Range("E33") "  +2.96%  "
Set-TextValue $ws.Range("D34") "3.76"
Set-TextValue $ws.Range("E34") "  +2.77%  "
Set-TextValue $ws.Range("D35") "25.13"
Set-TextValue $ws.Range("E35") "  +9.21%  "
Set-TextValue $ws.Range("E36") "  +10.81%  "
Set-TextValue $ws.Range("E37") "  +4.99%  "
Set-TextValue $ws.Range("E38") "  +5.89%  "
Set-TextValue $ws.Range("E39") "  +5.81%  "
Set-TextValue $ws.Range("E40") "  +18.30%  "
Set-TextValue $ws.Range("D41") "22.22"
Set-TextValue $ws.Range("E41") "  -0.34%  "
Set-TextValue $ws.Range("E42") "  -0.01%  "
Set-TextValue $ws.Range("D43") "0.415"
Set-TextValue $ws.Range("E43") "  +11.23%  "
Set-TextValue $ws.Range("E44") "  +17.83%  "
Set-TextValue $ws.Range("D45") "1.96"
Set-TextValue $ws.Range("E46") "  +0.02%  "
Set-TextValue $ws.Range("D47") "149.09"
Set-TextValue $ws.Range("E47") "  +2.09%  "
Set-TextValue $ws.Range("D48") "1.36"
Set-TextValue $ws.Range("E48") "  +10.55%  "
Set-TextValue $ws.Range("D49") "44.29"
Set-TextValue $ws.Range("E49") "  +1.69%  "
Set-TextValue $ws.Range("D50") "4.39"
Set-TextValue $ws.Range("E50") "  +9.07%  "
Set-TextValue $ws.Range("D51") "0.653"
Set-TextValue $ws.Range("E51") "  +9.70%  "
